$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the series. Insert a new row at
# position 154 (pushing the existing rows 154-174 down to 155-175) and
# populate it with the new observation's data.
$ws.Rows(154).Insert()

$ws.Range("A154").Value = 2
$ws.Range("B154").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C154").Value = "Coquimbo"
$ws.Range("D154").Value = 45015
$ws.Range("E154").Value = 4
$ws.Range("F154").Value = 100112024
$ws.Range("G154").Value = "Choclo"
$ws.Range("H154").Value = "Choclero"
$ws.Range("I154").Value = "Primera"
$ws.Range("J154").Value = 30000
$ws.Range("K154").Value = 250
$ws.Range("L154").Value = 300
$ws.Range("M154").Value = 275
$ws.Range("N154").Value = "$/unidad"
$ws.Range("O154").Value = "Provincia de Limarí"
$ws.Range("P154").Value = 275
$ws.Range("Q154").Value = 1
$ws.Range("R154").Value = "Hortaliza"
